$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.311.08"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.666.04"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.76%  "
$ws.Range("D5").Value = "219.40"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "0.5357"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "0.2661"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "0.06412"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").Value = "20.75"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "4.568"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "1.661.81"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "1.893.30"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "0.5538"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "0.0₅8247"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "65.94"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").Value = "194.48"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").Value = "10.30"
$ws.Range("E21").Value = "  +2.46%  "
$ws.Range("D22").Value = "6.044"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "146.24"
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("D25").Value = "0.1233"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "7.206"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "16.21"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "1.488"
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("D29").Value = "0.05838"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").Value = "3.618"
$ws.Range("E31").Value = "  +2.57%  "
$ws.Range("D32").Value = "3.282"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("D34").Value = "0.9714"
$ws.Range("E34").Value = "  +2.53%  "
$ws.Range("D35").Value = "2.826"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").Value = "2.420"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("D38").Value = "0.01609"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "0.8724"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.053.05"
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "105.20"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "1.804.62"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "57.92"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.014"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈104"
$ws.Range("E47").Value = "  -6.81%  "
$ws.Range("D48").Value = "0.4383"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").Value = "8.053"
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").Value = "0.05164"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "1.416"
